# Apply "cities update end October" edit
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the city name typo: "Bankok" -> "Bangkok" (row 38, column A)
$ws.Range("A38").Value = "Bangkok"

# Update counts (column B = in_deck, column C = in_reserve, column D = in_discard)
$ws.Range("B5").Value = 2

$ws.Range("B7").Value = 3
$ws.Range("C7").Value = 0

$ws.Range("B8").Value = 3
$ws.Range("C8").Value = 0

$ws.Range("B18").Value = 1
$ws.Range("C18").Value = 1

$ws.Range("B28").Value = 2

$ws.Range("B34").Value = 1

$ws.Range("B38").Value = 1
$ws.Range("B39").Value = 1
$ws.Range("B40").Value = 1
$ws.Range("B41").Value = 1
$ws.Range("B42").Value = 1
$ws.Range("B43").Value = 1
$ws.Range("B44").Value = 1
$ws.Range("B45").Value = 1

# Update the view: select B29 and scroll so row 4 is the top-left visible row
$ws.Range("B29").Select()
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1
